$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns C and D to fit the new Arabic text (closest achievable values
# given the engine's pixel-quantized column width storage: ~29.11 and ~65.82
# chars respectively).
$ws.Columns.Item(3).ColumnWidth = 28.25
$ws.Columns.Item(4).ColumnWidth = 64.92

# Arabic (ara) master-data rows, mirroring the existing eng/fra rows 2-7 / 8-13.
$rows = @(
    @{ Row=14; SrcRow=8;  Name="التسجيل المسبق"; Descr="بوابة الويب للتسجيلات المسبقة" },
    @{ Row=15; SrcRow=9;  Name="عميل التسجيل"; Descr="تطبيق سطح المكتب للتسجيلات" },
    @{ Row=16; SrcRow=10; Name="معالج التسجيل"; Descr="طلب عملية ما بعد التسجيل" },
    @{ Row=17; SrcRow=11; Name="مصادقة الهوية"; Descr="طلب مصادقة مزود خدمة الطرف الثالث" },
    @{ Row=18; SrcRow=12; Name="معرف التحكم"; Descr="بوابة الويب لتكوين التطبيقات" },
    @{ Row=19; SrcRow=13; Name="بوابة المقيمين"; Descr="بوابة الويب لخدمات إنشاء معرف البريد" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $srcRow = $r.SrcRow

    # Column A: lang_code = "ara"
    $ws.Cells.Item($rowNum, 1).Value = "ara"

    # Column B: numeric id, 10001-10006. The column has a Text ("@") style,
    # so flip to General while writing the literal number to keep it numeric.
    $ws.Cells.Item($rowNum, 2).NumberFormat = "General"
    $ws.Cells.Item($rowNum, 2).Value = $ws.Cells.Item($srcRow, 2).Value()
    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"

    # Column C: name (Arabic), wrapped + left aligned.
    $c = $ws.Cells.Item($rowNum, 3)
    $c.Value = $r.Name
    $c.HorizontalAlignment = -4131
    $c.WrapText = $true

    # Column D: descr (Arabic), wrapped + left aligned.
    $d = $ws.Cells.Item($rowNum, 4)
    $d.Value = $r.Descr
    $d.HorizontalAlignment = -4131
    $d.WrapText = $true

    # Column E: is_active = text "TRUE" (shared with existing rows). Typing the
    # literal word turns it into a native boolean, so copy the already-text
    # "TRUE" value from the matching French row instead.
    $ws.Range("E" + $srcRow).Copy()
    $ws.Range("E" + $rowNum).PasteSpecial(-4163)

    $ws.Rows.Item($rowNum).RowHeight = 16.4
}

# Move the active selection, matching the post-edit state recorded in the file.
$null = $ws.Range("D24").Select()
